$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.305.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.698.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +8.14%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.687.78"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.04%  "

# Row 8
$ws.Range("E8").Value = "  +4.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("E10").Value = "  +0.73%  "

# Row 11
$ws.Range("E11").Value = "  +26.33%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.611"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.06%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "49.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.29%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000288"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.88%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.296.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "678.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.40%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "9.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.93%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.698.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.447.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.81%  "

# Row 20
$ws.Range("E20").Value = "  +1.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.93%  "

# Row 22
$ws.Range("E22").Value = "  +2.41%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.943"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.66%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.36%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "102.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.63%  "

# Row 26
$ws.Range("E26").Value = "  +2.30%  "

# Row 27
$ws.Range("E27").Value = "  +6.86%  "

# Row 28
$ws.Range("E28").Value = "  +8.12%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.23%  "

# Row 31
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.53%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.25%  "

# Row 33
$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.05%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.69%  "

# Row 35
$ws.Range("B35").Value = "dogwifhat"
$ws.Range("C35").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.87%  "

# Row 36
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "585.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.93%  "

# Row 37
$ws.Range("B37").Value = "Cosmos"
$ws.Range("C37").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.91%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.109"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.52%  "

# Row 39
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.683.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.22%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.12%  "

# Row 42
$ws.Range("E42").Value = "  +5.15%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0458"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.31%  "

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.352"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.46%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0768"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.48%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.44%  "

# Row 47
$ws.Range("E47").Value = "  +4.24%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.38%  "

# Row 49
$ws.Range("E49").Value = "  +4.06%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.04%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.41%  "
